$wb = $excel.ActiveWorkbook

# ===== Sheet 1 (ALC) =====
$ws = $wb.Worksheets.Item(1)
$ws.Range("H19").Value = 1468.6154
$ws.Range("I19").Value = 1670.8889
$ws.Range("J19").Value = 1013.5
$ws.Range("K19").Value = 1670.8889
$ws.Range("L19").Value = 1013.5
$ws.Range("M19").Value = -1495.8889
$ws.Range("N19").Value = -1363.5
$ws.Range("H40").Value = 2545.919
$ws.Range("I40").Value = 2015
$ws.Range("J40").Value = 3526.077
$ws.Range("K40").Value = 2015
$ws.Range("L40").Value = 3526.077
$ws.Range("M40").Value = -1840
$ws.Range("N40").Value = -3876.077
$ws.Range("H43").Value = 5877.4443
$ws.Range("I43").Value = 5924.625
$ws.Range("K43").Value = 5924.625
$ws.Range("M43").Value = -5855.625
$ws.Range("H64").Value = 18185876
$ws.Range("I64").Value = 20003964
$ws.Range("K64").Value = 20003964
$ws.Range("M64").Value = -20003716
$ws.Range("H67").Value = 18185876
$ws.Range("I67").Value = 20003964
$ws.Range("K67").Value = 20003964
$ws.Range("M67").Value = -20003106
$ws.Range("H70").Value = 6391.1177
$ws.Range("I70").Value = 2492.7856
$ws.Range("K70").Value = 7478.3568
$ws.Range("M70").Value = -7208.3568
$ws.Range("H73").Value = 6391.1177
$ws.Range("I73").Value = 2492.7856
$ws.Range("K73").Value = 7478.3568
$ws.Range("M73").Value = -6542.3568
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 11073.207
$ws.Range("I86").Value = 10319.381
$ws.Range("K86").Value = 10319.381
$ws.Range("M86").Value = -9196.381
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 11073.207
$ws.Range("I89").Value = 10319.381
$ws.Range("K89").Value = 51596.905
$ws.Range("M89").Value = -45980.905
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H106").Value = 2817.6365
$ws.Range("J106").Value = 3128.2
$ws.Range("L106").Value = 3128.2
$ws.Range("N106").Value = -4390.2
$ws.Range("H125").Value = 1759.8
$ws.Range("I125").Value = 1954.6666
$ws.Range("K125").Value = 17591.9994
$ws.Range("M125").Value = -15131.9994
$ws.Range("H132").Value = 1633.3334
$ws.Range("I132").Value = 1726.4333
$ws.Range("K132").Value = 5179.2999
$ws.Range("M132").Value = -2649.2999
$ws.Range("H137").Value = 2476.3333
$ws.Range("I137").Value = 2423.7646
$ws.Range("J137").Value = 2699.75
$ws.Range("K137").Value = 7271.293799999999
$ws.Range("L137").Value = 8099.25
$ws.Range("M137").Value = -4721.293799999999
$ws.Range("N137").Value = -13199.25
$ws.Range("H141").Value = 4138.8
$ws.Range("I141").Value = 4866.3335
$ws.Range("J141").Value = 3047.5
$ws.Range("K141").Value = 14599.0005
$ws.Range("L141").Value = 9142.5
$ws.Range("M141").Value = -9419.000499999998
$ws.Range("N141").Value = -19502.5

# ===== Sheet 2 (ARM) =====
$ws = $wb.Worksheets.Item(2)
$ws.Range("H26").Value = 1217.3334
$ws.Range("I26").Value = 1217.3334
$ws.Range("K26").Value = 1217.3334
$ws.Range("M26").Value = -887.3334
$ws.Range("H32").Value = 4644.76
$ws.Range("I32").Value = 3310.3948
$ws.Range("K32").Value = 3310.3948
$ws.Range("M32").Value = -3023.3948
$ws.Range("H45").Value = 1651.25
$ws.Range("I45").Value = 1331
$ws.Range("J45").Value = 1971.5
$ws.Range("K45").Value = 1331
$ws.Range("L45").Value = 1971.5
$ws.Range("M45").Value = -954
$ws.Range("N45").Value = -2725.5
$ws.Range("H74").Value = 50005520
$ws.Range("I74").Value = 55561304
$ws.Range("K74").Value = 55561304
$ws.Range("M74").Value = -55560430
$ws.Range("H77").Value = 50005520
$ws.Range("I77").Value = 55561304
$ws.Range("K77").Value = 277806520
$ws.Range("M77").Value = -277802152
$ws.Range("H101").Value = 158336.67
$ws.Range("I101").Value = 200000
$ws.Range("J101").Value = 150004
$ws.Range("K101").Value = 200000
$ws.Range("L101").Value = 150004
$ws.Range("M101").Value = -196755
$ws.Range("N101").Value = -156494
$ws.Range("H102").Value = 20000950
$ws.Range("I102").Value = 50000450
$ws.Range("J102").Value = 1283.3334
$ws.Range("K102").Value = 50000450
$ws.Range("L102").Value = 1283.3334
$ws.Range("M102").Value = -49998828
$ws.Range("N102").Value = -4527.3334
$ws.Range("H110").Value = 44848.418
$ws.Range("I110").Value = 48471.5
$ws.Range("K110").Value = 48471.5
$ws.Range("M110").Value = -46426.5
$ws.Range("H132").Value = 3034139.8
$ws.Range("I132").Value = 3034139.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9102419.399999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9099889.399999999
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H135").Value = 99557.71
$ws.Range("J135").Value = 99557.71
$ws.Range("L135").Value = 99557.71
$ws.Range("N135").Value = -109697.71

# ===== Sheet 3 (BSM) =====
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 6899.2
$ws.Range("I20").Value = 4999
$ws.Range("K20").Value = 4999
$ws.Range("M20").Value = -4752
$ws.Range("H80").Value = 26316112
$ws.Range("J80").Value = 35714492
$ws.Range("L80").Value = 35714492
$ws.Range("N80").Value = -35716488
$ws.Range("H83").Value = 26316112
$ws.Range("J83").Value = 35714492
$ws.Range("L83").Value = 178572460
$ws.Range("N83").Value = -178582444
$ws.Range("H86").Value = 4668.615
$ws.Range("I86").Value = 3999.6
$ws.Range("K86").Value = 3999.6
$ws.Range("M86").Value = -2876.6
$ws.Range("H89").Value = 4668.615
$ws.Range("I89").Value = 3999.6
$ws.Range("K89").Value = 19998
$ws.Range("M89").Value = -14382
$ws.Range("H94").Value = 783.9
$ws.Range("I94").Value = 730.2857
$ws.Range("J94").Value = 909
$ws.Range("K94").Value = 730.2857
$ws.Range("L94").Value = 909
$ws.Range("M94").Value = -279.2857
$ws.Range("N94").Value = -1811
$ws.Range("H99").Value = 1875
$ws.Range("I99").Value = 1875
$ws.Range("K99").Value = 1875
$ws.Range("M99").Value = -377
$ws.Range("H105").Value = 2501.7778
$ws.Range("I105").Value = 2189.625
$ws.Range("K105").Value = 2189.625
$ws.Range("M105").Value = -442.625
$ws.Range("H134").Value = 71431650
$ws.Range("I134").Value = 100002300
$ws.Range("K134").Value = 300006900
$ws.Range("M134").Value = -300004365

# ===== Sheet 4 (CRP) =====
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 208.4
$ws.Range("J7").Value = 279.8
$ws.Range("L7").Value = 279.8
$ws.Range("N7").Value = -505.8
$ws.Range("H31").Value = 10740.108
$ws.Range("I31").Value = 7403.5
$ws.Range("K31").Value = 7403.5
$ws.Range("M31").Value = -7108.5
$ws.Range("H34").Value = 10740.108
$ws.Range("I34").Value = 7403.5
$ws.Range("K34").Value = 7403.5
$ws.Range("M34").Value = -7201.5
$ws.Range("H94").Value = 978.2
$ws.Range("I94").Value = 964.6667
$ws.Range("K94").Value = 964.6667
$ws.Range("M94").Value = -513.6667
$ws.Range("H99").Value = 1939.8182
$ws.Range("I99").Value = 1549.875
$ws.Range("J99").Value = 2979.6667
$ws.Range("K99").Value = 1549.875
$ws.Range("L99").Value = 2979.6667
$ws.Range("M99").Value = -51.875
$ws.Range("N99").Value = -5975.6667
$ws.Range("H103").Value = 35397.8
$ws.Range("I103").Value = 26747.5
$ws.Range("K103").Value = 26747.5
$ws.Range("M103").Value = -25575.5
$ws.Range("H122").Value = 4323
$ws.Range("I122").Value = 4718.5
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 14155.5
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -11705.5
$ws.Range("N122").Value = -10750
$ws.Range("H126").Value = 1939.8182
$ws.Range("I126").Value = 1549.875
$ws.Range("J126").Value = 2979.6667
$ws.Range("K126").Value = 4649.625
$ws.Range("L126").Value = 8939.000100000001
$ws.Range("M126").Value = -2179.625
$ws.Range("N126").Value = -13879.0001
$ws.Range("H132").Value = 31402168
$ws.Range("I132").Value = 34392790
$ws.Range("J132").Value = 621.25
$ws.Range("K132").Value = 103178370
$ws.Range("L132").Value = 1863.75
$ws.Range("M132").Value = -103175840
$ws.Range("N132").Value = -6923.75
$ws.Range("H134").Value = 25001858
$ws.Range("I134").Value = 41667940
$ws.Range("J134").Value = 2734.75
$ws.Range("K134").Value = 125003820
$ws.Range("L134").Value = 8204.25
$ws.Range("M134").Value = -125001285
$ws.Range("N134").Value = -13274.25

# ===== Sheet 5 (CUL) =====
$ws = $wb.Worksheets.Item(5)
$ws.Range("H70").Value = 10931
$ws.Range("I70").Value = 6399
$ws.Range("J70").Value = 19995
$ws.Range("K70").Value = 19197
$ws.Range("L70").Value = 59985
$ws.Range("M70").Value = -18882
$ws.Range("N70").Value = -60615
$ws.Range("H73").Value = 10931
$ws.Range("I73").Value = 6399
$ws.Range("J73").Value = 19995
$ws.Range("K73").Value = 19197
$ws.Range("L73").Value = 59985
$ws.Range("M73").Value = -18105
$ws.Range("N73").Value = -62169
$ws.Range("H120").Value = 29998.5
$ws.Range("I120").Value = 29998
$ws.Range("K120").Value = 89994
$ws.Range("M120").Value = -85156
$ws.Range("H121").Value = 65270.26
$ws.Range("I121").Value = 102172.3
$ws.Range("K121").Value = 306516.9
$ws.Range("M121").Value = -305206.9
$ws.Range("H123").Value = 7959
$ws.Range("J123").Value = 7959
$ws.Range("L123").Value = 23877
$ws.Range("N123").Value = -28777

# ===== Sheet 6 (GSM) =====
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 1335.7778
$ws.Range("I2").Value = 1800.3334
$ws.Range("J2").Value = 406.66666
$ws.Range("K2").Value = 1800.3334
$ws.Range("L2").Value = 406.66666
$ws.Range("M2").Value = -1687.3334
$ws.Range("N2").Value = -632.66666
$ws.Range("H3").Value = 900.75
$ws.Range("J3").Value = 867.6667
$ws.Range("L3").Value = 867.6667
$ws.Range("N3").Value = -1099.6667
$ws.Range("H7").Value = 929088.2
$ws.Range("J7").Value = 21710.857
$ws.Range("L7").Value = 21710.857
$ws.Range("N7").Value = -21934.857
$ws.Range("H8").Value = 929088.2
$ws.Range("J8").Value = 21710.857
$ws.Range("L8").Value = 21710.857
$ws.Range("N8").Value = -21988.857
$ws.Range("H10").Value = 6750.75
$ws.Range("I10").Value = 7667.6665
$ws.Range("K10").Value = 7667.6665
$ws.Range("M10").Value = -7498.6665
$ws.Range("H11").Value = 37513376
$ws.Range("I11").Value = 50017668
$ws.Range("K11").Value = 50017668
$ws.Range("M11").Value = -50017529
$ws.Range("H14").Value = 2036034.6
$ws.Range("J14").Value = 60001.332
$ws.Range("L14").Value = 60001.332
$ws.Range("N14").Value = -60337.332
$ws.Range("H20").Value = 3343366.2
$ws.Range("I20").Value = 3343366.2
$ws.Range("K20").Value = 3343366.2
$ws.Range("M20").Value = -3343121.2
$ws.Range("H21").Value = 10
$ws.Range("I21").Value = 10
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 10
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 163
$ws.Range("N21").ClearContents()
$ws.Range("H24").Value = 37999
$ws.Range("I24").Value = 37999
$ws.Range("K24").Value = 37999
$ws.Range("M24").Value = -37826
$ws.Range("H29").Value = 6936.6665
$ws.Range("I29").Value = 6936.6665
$ws.Range("K29").Value = 6936.6665
$ws.Range("M29").Value = -6646.6665
$ws.Range("H30").Value = 10
$ws.Range("I30").Value = 10
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 10
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 95
$ws.Range("N30").ClearContents()
$ws.Range("H108").Value = 76500
$ws.Range("J108").Value = 76500
$ws.Range("L108").Value = 76500
$ws.Range("N108").Value = -84180
$ws.Range("H113").Value = 94495.09
$ws.Range("I113").Value = 113494.11
$ws.Range("J113").Value = 8999.5
$ws.Range("K113").Value = 113494.11
$ws.Range("L113").Value = 8999.5
$ws.Range("M113").Value = -111324.11
$ws.Range("N113").Value = -13339.5
$ws.Range("H126").Value = 3650.6316
$ws.Range("I126").Value = 4091.5334
$ws.Range("J126").Value = 1997.25
$ws.Range("K126").Value = 12274.6002
$ws.Range("L126").Value = 5991.75
$ws.Range("M126").Value = -9804.600199999999
$ws.Range("N126").Value = -10931.75
$ws.Range("H132").Value = 31257780
$ws.Range("I132").Value = 41675188
$ws.Range("J132").Value = 5554
$ws.Range("K132").Value = 125025564
$ws.Range("L132").Value = 16662
$ws.Range("M132").Value = -125023034
$ws.Range("N132").Value = -21722

# ===== Sheet 7 (LTW) =====
$ws = $wb.Worksheets.Item(7)
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H20").Value = 45000
$ws.Range("J20").Value = 45000
$ws.Range("L20").Value = 45000
$ws.Range("N20").Value = -45452
$ws.Range("H22").Value = 4118.5713
$ws.Range("I22").Value = 3270
$ws.Range("J22").Value = 5250
$ws.Range("K22").Value = 3270
$ws.Range("L22").Value = 5250
$ws.Range("M22").Value = -2975
$ws.Range("N22").Value = -5840
$ws.Range("H27").Value = 4118.5713
$ws.Range("I27").Value = 3270
$ws.Range("J27").Value = 5250
$ws.Range("K27").Value = 3270
$ws.Range("L27").Value = 5250
$ws.Range("M27").Value = -3163
$ws.Range("N27").Value = -5464
$ws.Range("H40").Value = 2011.65
$ws.Range("I40").Value = 1912.8889
$ws.Range("K40").Value = 1912.8889
$ws.Range("M40").Value = -1776.8889
$ws.Range("H46").Value = 1300
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 1400
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 1400
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -1776
$ws.Range("H55").Value = 225.3158
$ws.Range("I55").Value = 198.9
$ws.Range("J55").Value = 254.66667
$ws.Range("K55").Value = 198.9
$ws.Range("L55").Value = 254.66667
$ws.Range("M55").Value = -25.90000000000001
$ws.Range("N55").Value = -600.6666700000001
$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 30000
$ws.Range("K69").Value = 30000
$ws.Range("M69").Value = -29189
$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 30000
$ws.Range("K72").Value = 90000
$ws.Range("M72").Value = -85944
$ws.Range("H93").Value = 1370.2609
$ws.Range("I93").Value = 1410.1177
$ws.Range("J93").Value = 1257.3334
$ws.Range("K93").Value = 1410.1177
$ws.Range("L93").Value = 1257.3334
$ws.Range("M93").Value = -162.1177
$ws.Range("N93").Value = -3753.3334
$ws.Range("H130").Value = 93475.664
$ws.Range("J130").Value = 93475.664
$ws.Range("L130").Value = 93475.664
$ws.Range("N130").Value = -103515.664
$ws.Range("H132").Value = 17873878
$ws.Range("I132").Value = 19248614
$ws.Range("K132").Value = 57745842
$ws.Range("M132").Value = -57743312
$ws.Range("H134").Value = 81999
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 81999
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 81999
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -92139

# ===== Sheet 8 (WVR) =====
$ws = $wb.Worksheets.Item(8)
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H52").Value = 14800
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880
$ws.Range("H107").Value = 1490.625
$ws.Range("I107").Value = 1380
$ws.Range("J107").Value = 1576.6666
$ws.Range("K107").Value = 4140
$ws.Range("L107").Value = 4729.9998
$ws.Range("M107").Value = -2220
$ws.Range("N107").Value = -8569.9998
$ws.Range("H113").Value = 871.3333
$ws.Range("I113").Value = 871.3333
$ws.Range("K113").Value = 2613.9999
$ws.Range("M113").Value = -443.9998999999998
$ws.Range("H122").Value = 4611.7646
$ws.Range("I122").Value = 4874.8
$ws.Range("K122").Value = 14624.4
$ws.Range("M122").Value = -12174.4
$ws.Range("H132").Value = 14711496
$ws.Range("I132").Value = 26318666
$ws.Range("K132").Value = 78955998
$ws.Range("M132").Value = -78953468
